$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68").Value = "GRT-USD"
